# Generate Report for Handoff
# The localization file 61cc4c18-af55-47a0-9d1f-f05319221a94 is now ready
# for handoff: update its Status / Priority / Latest Handoff Datetime on
# the per-locale sheets, and mirror the Status + generate-date on the
# Overview roll-up sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row 3 (61cc4c18-af55-47a0-9d1f-f05319221a94.md)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-24 00:14:31"

# ---------------------------------------------------------------------
# zh-cn sheet - row 3 (61cc4c18-af55-47a0-9d1f-f05319221a94.md)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-24 00:14:27"

# ---------------------------------------------------------------------
# de-de sheet - row 3 (61cc4c18-af55-47a0-9d1f-f05319221a94.md)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-24 00:14:31"

# ---------------------------------------------------------------------
# Column widths widened to fit the new, longer "Ready for handoff" text
# (closest width this runtime's ColumnWidth rounding can reach to the
# recorded 17.2159881591797 target).
# ---------------------------------------------------------------------
$wsOverview.Range("E:F").ColumnWidth = 16.3
$wsZhCn.Range("C:C").ColumnWidth = 16.3
$wsDeDe.Range("C:C").ColumnWidth = 16.3
